$wb = $excel.ActiveWorkbook

# 1) Set the selection on SARIMA (currently active/tabSelected sheet) to G26
#    before we switch the active sheet away from it, so SARIMA keeps a
#    real selection instead of defaulting back to A1.
$sarima = $wb.Worksheets.Item("SARIMA")
$sarima.Activate()
$sarima.Range("G26").Select()

# 2) Add the new "Stage" worksheet as the LAST tab.
$count = $wb.Worksheets.Count
$stage = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($count))
$stage.Name = "Stage"

# 3) Populate header + data. Column A first top-to-bottom, then column B,
#    leaving B1 ("Stage") for last so the shared-string table grows in the
#    same order as the source workbook (Industry, Manufacturer, Retailer,
#    Wholesaler, Stage).
$stage.Range("A1").Value = "Industry"
$stage.Range("A2").Value = "Apparel"
$stage.Range("B2").Value = "Manufacturer"
$stage.Range("A3").Value = "Audio and Video Equipment"
$stage.Range("B3").Value = "Manufacturer"
$stage.Range("A4").Value = "Automobile Manufacturing"
$stage.Range("B4").Value = "Manufacturer"
$stage.Range("A5").Value = "Battery Manufacturing"
$stage.Range("B5").Value = "Manufacturer"
$stage.Range("A6").Value = "Beverage Manufacturing"
$stage.Range("B6").Value = "Manufacturer"
$stage.Range("A7").Value = "Communications Equipment Manufacturing, Defense"
$stage.Range("B7").Value = "Manufacturer"
$stage.Range("A8").Value = "Communications Equipment Manufacturing, Nondefense"
$stage.Range("B8").Value = "Manufacturer"
$stage.Range("A9").Value = "Computer Storage Device Manufacturing"
$stage.Range("B9").Value = "Manufacturer"
$stage.Range("A10").Value = "Construction Machinery Manufacturing"
$stage.Range("B10").Value = "Manufacturer"
$stage.Range("A11").Value = "Dairy Product Manufacturing"
$stage.Range("B11").Value = "Manufacturer"
$stage.Range("A12").Value = "Electric Lighting Equipment Manufacturing"
$stage.Range("B12").Value = "Manufacturer"
$stage.Range("A13").Value = "Electrical Equipment Manufacturing"
$stage.Range("B13").Value = "Manufacturer"
$stage.Range("A14").Value = "Electromedical, Measuring, and Control Instrument Manufacturing"
$stage.Range("B14").Value = "Manufacturer"
$stage.Range("A15").Value = "Electronic Computer Manufacturing"
$stage.Range("B15").Value = "Manufacturer"
$stage.Range("A16").Value = "Fabricated Metal Products"
$stage.Range("B16").Value = "Manufacturer"
$stage.Range("A17").Value = "Farm Machinery and Equipment Manufacturing"
$stage.Range("B17").Value = "Manufacturer"
$stage.Range("A18").Value = "Ferrous Metal Foundries"
$stage.Range("B18").Value = "Manufacturer"
$stage.Range("A19").Value = "Furniture and Related Products"
$stage.Range("B19").Value = "Manufacturer"
$stage.Range("A20").Value = "Grain and Oilseed Milling"
$stage.Range("B20").Value = "Manufacturer"
$stage.Range("A21").Value = "Heavy Duty Truck Manufacturing"
$stage.Range("B21").Value = "Manufacturer"
$stage.Range("A22").Value = "Household Appliance Manufacturing"
$stage.Range("B22").Value = "Manufacturer"
$stage.Range("A23").Value = "Industrial Machinery Manufacturing"
$stage.Range("B23").Value = "Manufacturer"
$stage.Range("A24").Value = "Iron and Steel Mills and Ferroalloy and Steel Product Manufacturing"
$stage.Range("B24").Value = "Manufacturer"
$stage.Range("A25").Value = "Leather and Allied Products"
$stage.Range("B25").Value = "Manufacturer"
$stage.Range("A26").Value = "Light Truck and Utility Vehicle Manufacturing"
$stage.Range("B26").Value = "Manufacturer"
$stage.Range("A27").Value = "Material Handling Equipment Manufacturing"
$stage.Range("B27").Value = "Manufacturer"
$stage.Range("A28").Value = "Meat, Poultry, and Seafood Product Processing"
$stage.Range("B28").Value = "Manufacturer"
$stage.Range("A29").Value = "Metalworking Machinery Manufacturing"
$stage.Range("B29").Value = "Manufacturer"
$stage.Range("A30").Value = "Mining, Oil, and Gas Field Machinery Manufacturing"
$stage.Range("B30").Value = "Manufacturer"
$stage.Range("A31").Value = "Miscellaneous Products"
$stage.Range("B31").Value = "Manufacturer"
$stage.Range("A32").Value = "Nonmetallic Mineral Products"
$stage.Range("B32").Value = "Manufacturer"
$stage.Range("A33").Value = "Other Computer Peripheral Equipment Manufacturing"
$stage.Range("B33").Value = "Manufacturer"
$stage.Range("A34").Value = "Other Electronic Component Manufacturing"
$stage.Range("B34").Value = "Manufacturer"
$stage.Range("A35").Value = "Paint, Coating, and Adhesive, Manufacturing"
$stage.Range("B35").Value = "Manufacturer"
$stage.Range("A36").Value = "Paperboard Container Manufacturing"
$stage.Range("B36").Value = "Manufacturer"
$stage.Range("A37").Value = "Pesticide, Fertilizer, and Other Agricultural Chemical Manufacturing"
$stage.Range("B37").Value = "Manufacturer"
$stage.Range("A38").Value = "Petroleum and Coal Products"
$stage.Range("B38").Value = "Manufacturer"
$stage.Range("A39").Value = "Pharmaceutical and Medicine Manufacturing"
$stage.Range("B39").Value = "Manufacturer"
$stage.Range("A40").Value = "Photographic Equipment Manufacturing"
$stage.Range("B40").Value = "Manufacturer"
$stage.Range("A41").Value = "Plastics and Rubber Products"
$stage.Range("B41").Value = "Manufacturer"
$stage.Range("A42").Value = "Printing"
$stage.Range("B42").Value = "Manufacturer"
$stage.Range("A43").Value = "Pulp, Paper, and Paperboard Mills"
$stage.Range("B43").Value = "Manufacturer"
$stage.Range("A44").Value = "Search and Navigation Equipment, Defense"
$stage.Range("B44").Value = "Manufacturer"
$stage.Range("A45").Value = "Search and Navigation Equipment, Nondefense"
$stage.Range("B45").Value = "Manufacturer"
$stage.Range("A46").Value = "Textile Mills"
$stage.Range("B46").Value = "Manufacturer"
$stage.Range("A47").Value = "Textile Products"
$stage.Range("B47").Value = "Manufacturer"
$stage.Range("A48").Value = "Tobacco Manufacturing"
$stage.Range("B48").Value = "Manufacturer"
$stage.Range("A49").Value = "Total Manufacturing"
$stage.Range("B49").Value = "Manufacturer"
$stage.Range("A50").Value = "Transportation Equipment"
$stage.Range("B50").Value = "Manufacturer"
$stage.Range("A51").Value = "Ventilation, Heating, Air-conditioning, and Refrigeration Equipment Manufacturing"
$stage.Range("B51").Value = "Manufacturer"
$stage.Range("A52").Value = "Wood Products"
$stage.Range("B52").Value = "Manufacturer"
$stage.Range("A53").Value = "Building Mat. and Garden Equip. and Supplies Dealers"
$stage.Range("B53").Value = "Retailer"
$stage.Range("A54").Value = "Clothing and Clothing Access. Stores"
$stage.Range("B54").Value = "Retailer"
$stage.Range("A55").Value = "Food and Beverage Stores"
$stage.Range("B55").Value = "Retailer"
$stage.Range("A56").Value = "Furniture, Home Furn, Electronics, and Appliance Stores"
$stage.Range("B56").Value = "Retailer"
$stage.Range("A57").Value = "General Merchandise Stores"
$stage.Range("B57").Value = "Retailer"
$stage.Range("A58").Value = "Motor Vehicle and Parts Dealers"
$stage.Range("B58").Value = "Retailer"
$stage.Range("A59").Value = "Retail Trade"
$stage.Range("B59").Value = "Retailer"
$stage.Range("A60").Value = "Apparel, Piece Goods, and Notions"
$stage.Range("B60").Value = "Wholesaler"
$stage.Range("A61").Value = "Beer, Wine, and Distilled Alcoholic Beverages"
$stage.Range("B61").Value = "Wholesaler"
$stage.Range("A62").Value = "Chemicals and Allied Products"
$stage.Range("B62").Value = "Wholesaler"
$stage.Range("A63").Value = "Drugs and Druggists' Sundries"
$stage.Range("B63").Value = "Wholesaler"
$stage.Range("A64").Value = "Farm Product Raw Materials"
$stage.Range("B64").Value = "Wholesaler"
$stage.Range("A65").Value = "Furniture and Home Furnishings"
$stage.Range("B65").Value = "Wholesaler"
$stage.Range("A66").Value = "Grocery and Related Products"
$stage.Range("B66").Value = "Wholesaler"
$stage.Range("A67").Value = "Hardware, and Plumbing and Heating Equipment and Supplies"
$stage.Range("B67").Value = "Wholesaler"
$stage.Range("A68").Value = "Household Appliances and Electrical and Electronic Goods"
$stage.Range("B68").Value = "Wholesaler"
$stage.Range("A69").Value = "Lumber and Other Construction Materials"
$stage.Range("B69").Value = "Wholesaler"
$stage.Range("A70").Value = "Machinery, Equipment, and Supplies"
$stage.Range("B70").Value = "Wholesaler"
$stage.Range("A71").Value = "Metals and Minerals, Except Petroleum"
$stage.Range("B71").Value = "Wholesaler"
$stage.Range("A72").Value = "Miscellaneous Durable Goods"
$stage.Range("B72").Value = "Wholesaler"
$stage.Range("A73").Value = "Miscellaneous Nondurable Goods"
$stage.Range("B73").Value = "Wholesaler"
$stage.Range("A74").Value = "Motor Vehicle and Motor Vehicle Parts and Supplies"
$stage.Range("B74").Value = "Wholesaler"
$stage.Range("A75").Value = "Paper and Paper Products"
$stage.Range("B75").Value = "Wholesaler"
$stage.Range("A76").Value = "Petroleum and Petroleum Products"
$stage.Range("B76").Value = "Wholesaler"
$stage.Range("A77").Value = "Professional and Commercial Equipment and Supplies"
$stage.Range("B77").Value = "Wholesaler"
$stage.Range("A78").Value = "Total Merchant Wholesalers, Except Manufacturers' Sales Branches and Offices"
$stage.Range("B78").Value = "Wholesaler"
$stage.Range("B1").Value = "Stage"

# 4) Column widths (best-fit in the source file).
$stage.Columns.Item(1).ColumnWidth = 64.5
$stage.Columns.Item(2).ColumnWidth = 10.666666666666666

# 5) Selection + active sheet: Stage ends up active/selected at A10.
$stage.Range("A10").Select()
$stage.Activate()

Write-Output "Sheets: $($wb.Worksheets | ForEach-Object { $_.Name })"
